# Update image loading method.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Append note about the faster article table loading method.
$ws.Range("C39").Value = "Load multiply images on the screen. Add progress indicator. Add swipe right animation. Update UI, and delete unused code. Update article table loading method. It's much faster now."

# 2. Add an extra hour of work to the formula for that day.
$ws.Range("B39").Formula = "=2+0.5+0.5+1.5+1"
